$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated "Price" values are plain numeric-looking strings (e.g. "215.04").
# The source data stores these as text (inlineStr), so force a Text number format
# before assigning, otherwise Excel would auto-convert them to actual numbers.
$textCells = @(
    "D5", "D10", "D14", "D15", "D16", "D18", "D19", "D22", "D24", "D25", "D26", "D34", "D39", "D40", "D41", "D42", "D43", "D44", "D47", "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values (Coin name / Link / Price / Volume(1h))
$ws.Range("D2").Value = "26.972.61"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "1.679.12"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "215.04"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("E6").Value = "  -3.28%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").Value = "20.38"
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").Value = "1.915.57"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").Value = "1.688.76"
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("D14").Value = "4.10"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").Value = "0.528"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "65.72"
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "26.987.92"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "8.17"
$ws.Range("E18").Value = "  +6.11%  "
$ws.Range("D19").Value = "236.57"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("D20").Value = "0.0₃0734"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "4.45"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("D24").Value = "2.15"
$ws.Range("E24").Value = "  -3.03%  "
$ws.Range("D25").Value = "146.65"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("D26").Value = "7.21"
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("E27").Value = "  +1.37%  "
$ws.Range("E28").Value = "  -4.18%  "
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").Value = "1.483.60"
$ws.Range("E33").Value = "  +1.89%  "
$ws.Range("D34").Value = "3.17"
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("E35").Value = "  +5.38%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("E37").Value = "  +3.62%  "
$ws.Range("E38").Value = "  +3.85%  "
$ws.Range("D39").Value = "0.906"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("D40").Value = "5.80"
$ws.Range("E40").Value = "  -3.84%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "1.02"
$ws.Range("E41").Value = "  +3.29%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").Value = "2.32"
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("D44").Value = "67.30"
$ws.Range("E44").Value = "  +2.67%  "
$ws.Range("D45").Value = "1.819.21"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").Value = "90.64"
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("E48").Value = "  +4.29%  "
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").Value = "0.103"
$ws.Range("E50").Value = "  +2.99%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.75"
$ws.Range("E51").Value = "  +1.90%  "
